# Add two new Mac-Address rows (31 and 32) to the worksheet, matching the
# pattern of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data: regcntr_id, usr_id, machine_id, lang_code, is_active, cr_by, cr_dtimes, eff_dtimes
$newRows = @(
    @(10001, 110030, 10030),
    @(10001, 110031, 10031)
)

$startRow = 31
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
}

# Update the selection to match the saved workbook state (Excel records the
# last active cell in the sheet view when the file is saved).
$ws.Range("F30").Select()
